# Updated cryptos list (prices / 1h volume) per GitHub Actions refresh.
# For D-column cells whose new text looks like a plain number (e.g. "22.57"),
# we briefly force a Text number format so Excel stores the literal string
# instead of silently converting it to a numeric value, then restore the
# cell's original style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.910.60'
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').Value = '1.970.68'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('E4').Value = '  +0.08%  '
$orig = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '326.37'
$ws.Range('D5').Style = $orig
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('E6').Value = '  +0.12%  '
$orig = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4821'
$ws.Range('D7').Style = $orig
$ws.Range('E7').Value = '  -4.46%  '
$orig = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4080'
$ws.Range('D8').Style = $orig
$ws.Range('E8').Value = '  -4.43%  '
$orig = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '53.91'
$ws.Range('D9').Style = $orig
$ws.Range('E9').Value = '  -1.79%  '
$orig = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08647'
$ws.Range('D10').Style = $orig
$ws.Range('E10').Value = '  -6.16%  '
$orig = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.066'
$ws.Range('D11').Style = $orig
$ws.Range('E11').Value = '  -5.55%  '
$orig = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.57'
$ws.Range('D12').Style = $orig
$ws.Range('E12').Value = '  -4.39%  '
$ws.Range('D13').Value = '1.951.09'
$ws.Range('E13').Value = '  -1.69%  '
$orig = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.758'
$ws.Range('D14').Style = $orig
$ws.Range('E14').Value = '  -4.95%  '
$orig = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.287'
$ws.Range('D15').Style = $orig
$ws.Range('E15').Value = '  -4.16%  '
$orig = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.014'
$ws.Range('D16').Style = $orig
$ws.Range('E16').Value = '  +0.14%  '
$orig = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '91.01'
$ws.Range('D17').Style = $orig
$ws.Range('E17').Value = '  -4.74%  '
$orig = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001075'
$ws.Range('D18').Style = $orig
$ws.Range('E18').Value = '  -4.65%  '
$orig = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06634'
$ws.Range('D19').Style = $orig
$ws.Range('E19').Value = '  -0.72%  '
$orig = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.87'
$ws.Range('D20').Style = $orig
$ws.Range('E20').Value = '  -5.43%  '
$ws.Range('E21').Value = '  +0.17%  '
$orig = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.820'
$ws.Range('D22').Style = $orig
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').Value = '28.927.21'
$ws.Range('E23').Value = '  -2.41%  '
$orig = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.63'
$ws.Range('D24').Style = $orig
$ws.Range('E24').Value = '  -3.86%  '
$orig = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.290'
$ws.Range('D25').Style = $orig
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('D26').Value = '2.187.48'
$ws.Range('E26').Value = '  -1.81%  '
$orig = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.41'
$ws.Range('D27').Style = $orig
$ws.Range('E27').Value = '  -2.13%  '
$orig = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '154.25'
$ws.Range('D28').Style = $orig
$ws.Range('E28').Value = '  -3.20%  '
$orig = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.073'
$ws.Range('D29').Style = $orig
$ws.Range('E29').Value = '  -5.78%  '
$orig = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.181'
$ws.Range('D30').Style = $orig
$ws.Range('E30').Value = '  -6.71%  '
$orig = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '125.14'
$ws.Range('D31').Style = $orig
$ws.Range('E31').Value = '  -3.03%  '
$orig = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.015'
$ws.Range('D32').Style = $orig
$ws.Range('E32').Value = '  -5.06%  '
$orig = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09678'
$ws.Range('D33').Style = $orig
$ws.Range('E33').Value = '  -2.97%  '
$orig = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.478'
$ws.Range('D34').Style = $orig
$ws.Range('E34').Value = '  -6.87%  '
$orig = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.737'
$ws.Range('D35').Style = $orig
$ws.Range('E35').Value = '  -2.34%  '
$orig = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.702'
$ws.Range('D36').Style = $orig
$ws.Range('E36').Value = '  -2.96%  '
$orig = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02365'
$ws.Range('D37').Style = $orig
$ws.Range('E37').Value = '  -4.76%  '
$orig = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.284'
$ws.Range('D38').Style = $orig
$ws.Range('E38').Value = '  -2.90%  '
$orig = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06283'
$ws.Range('D39').Style = $orig
$ws.Range('E39').Value = '  -1.78%  '
$orig = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.862'
$ws.Range('D40').Style = $orig
$ws.Range('E40').Value = '  -7.52%  '
$orig = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6295'
$ws.Range('D41').Style = $orig
$ws.Range('E41').Value = '  -4.81%  '
$orig = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.25'
$ws.Range('D42').Style = $orig
$ws.Range('E42').Value = '  -4.78%  '
$orig = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.011'
$ws.Range('D43').Style = $orig
$ws.Range('E43').Value = '  +0.11%  '
$orig = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1929'
$ws.Range('D44').Style = $orig
$ws.Range('E44').Value = '  -7.14%  '
$orig = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.334'
$ws.Range('D45').Style = $orig
$ws.Range('E45').Value = '  +2.77%  '
$orig = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6012'
$ws.Range('D46').Style = $orig
$ws.Range('E46').Value = '  -5.85%  '
$ws.Range('E47').Value = '  -3.59%  '
$orig = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.103'
$ws.Range('D48').Style = $orig
$ws.Range('E48').Value = '  -5.43%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$orig = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00000000340'
$ws.Range('D49').Style = $orig
$ws.Range('E49').Value = '  +5.18%  '
$ws.Range('B50').Value = 'PancakeSwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$orig = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.449'
$ws.Range('D50').Style = $orig
$ws.Range('E50').Value = '  -2.38%  '
$orig = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.147'
$ws.Range('D51').Style = $orig
$ws.Range('E51').Value = '  +6.67%  '
